# T1753 test data: rename/re-point the "Repeaters" tabs and update the
# active sheet / selection state to match.
#
# Before: Repeaters (old data), Repeaters_Updated (new data, was active,
#         cell B16 selected)
# After:  RepeatersOld (old data, renamed), Repeaters (was
#         Repeaters_Updated, now active with A17 selected)

$wb = $excel.ActiveWorkbook

# Rename the current "Repeaters" sheet out of the way first so the name
# "Repeaters" is free for the renamed "Repeaters_Updated" sheet.
$repeatersOld = $wb.Worksheets.Item("Repeaters")
$repeatersOld.Name = "RepeatersOld"

$repeatersNew = $wb.Worksheets.Item("Repeaters_Updated")
$repeatersNew.Name = "Repeaters"

# Make the renamed "Repeaters" sheet the active tab with A17 selected.
$repeatersNew.Activate()
$repeatersNew.Range("A17").Select()
